$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 43.62448933333334
$ws.Cells.Item(2, 8).Value = 130.873468
$ws.Cells.Item(2, 9).Value = 0.289568119079398
$ws.Cells.Item(2, 10).Value = 0.2895681190793979
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 19.00842133333333
$ws.Cells.Item(2, 14).Value = 57.02526399999999
$ws.Cells.Item(2, 15).Value = 0.2688244872258518
$ws.Cells.Item(2, 16).Value = 0.2688244872258518
$ws.Cells.Item(2, 17).Value = 829.2326736995058
$ws.Cells.Item(2, 18).Value = 7463.094063295551
$ws.Cells.Item(2, 19).Value = 0.07784300112847357
$ws.Cells.Item(2, 20).Value = 0.07784300112847355

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 43.62448933333334
$ws.Cells.Item(3, 8).Value = 130.873468
$ws.Cells.Item(3, 9).Value = 0.289568119079398
$ws.Cells.Item(3, 10).Value = 0.2895681190793979
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 21.993491
$ws.Cells.Item(3, 14).Value = 65.980473
$ws.Cells.Item(3, 15).Value = 0.3110405034011621
$ws.Cells.Item(3, 16).Value = 0.311040503401162
$ws.Cells.Item(3, 17).Value = 959.4548135322628
$ws.Cells.Item(3, 18).Value = 8635.093321790364
$ws.Cells.Item(3, 19).Value = 0.09006741352738359
$ws.Cells.Item(3, 20).Value = 0.09006741352738357

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 43.62448933333334
$ws.Cells.Item(4, 8).Value = 130.873468
$ws.Cells.Item(4, 9).Value = 0.289568119079398
$ws.Cells.Item(4, 10).Value = 0.2895681190793979
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 17.47796633333333
$ws.Cells.Item(4, 14).Value = 52.433899
$ws.Cells.Item(4, 15).Value = 0.2471801973933361
$ws.Cells.Item(4, 16).Value = 0.2471801973933361
$ws.Cells.Item(4, 17).Value = 762.4673558768591
$ws.Cells.Item(4, 18).Value = 6862.206202891732
$ws.Cells.Item(4, 19).Value = 0.07157550483286267
$ws.Cells.Item(4, 20).Value = 0.07157550483286265

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 43.62448933333334
$ws.Cells.Item(5, 8).Value = 130.873468
$ws.Cells.Item(5, 9).Value = 0.289568119079398
$ws.Cells.Item(5, 10).Value = 0.2895681190793979
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 12.229533
$ws.Cells.Item(5, 14).Value = 36.688599
$ws.Cells.Item(5, 15).Value = 0.17295481197965
$ws.Cells.Item(5, 16).Value = 0.17295481197965
$ws.Cells.Item(5, 17).Value = 533.5071319101481
$ws.Cells.Item(5, 18).Value = 4801.564187191332
$ws.Cells.Item(5, 19).Value = 0.05008219959067818
$ws.Cells.Item(5, 20).Value = 0.05008219959067816

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 39.73808533333334
$ws.Cells.Item(6, 8).Value = 119.214256
$ws.Cells.Item(6, 9).Value = 0.2637711707721372
$ws.Cells.Item(6, 10).Value = 0.2637711707721372
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 19.00842133333333
$ws.Cells.Item(6, 14).Value = 57.02526399999999
$ws.Cells.Item(6, 15).Value = 0.2688244872258518
$ws.Cells.Item(6, 16).Value = 0.2688244872258518
$ws.Cells.Item(6, 17).Value = 755.3582689959537
$ws.Cells.Item(6, 18).Value = 6798.224420963584
$ws.Cells.Item(6, 19).Value = 0.07090814972778239
$ws.Cells.Item(6, 20).Value = 0.07090814972778237

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 39.73808533333334
$ws.Cells.Item(7, 8).Value = 119.214256
$ws.Cells.Item(7, 9).Value = 0.2637711707721372
$ws.Cells.Item(7, 10).Value = 0.2637711707721372
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 21.993491
$ws.Cells.Item(7, 14).Value = 65.980473
$ws.Cells.Item(7, 15).Value = 0.3110405034011621
$ws.Cells.Item(7, 16).Value = 0.311040503401162
$ws.Cells.Item(7, 17).Value = 873.9792221358989
$ws.Cells.Item(7, 18).Value = 7865.812999223089
$ws.Cells.Item(7, 19).Value = 0.08204351773967945
$ws.Cells.Item(7, 20).Value = 0.08204351773967943

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 39.73808533333334
$ws.Cells.Item(8, 8).Value = 119.214256
$ws.Cells.Item(8, 9).Value = 0.2637711707721372
$ws.Cells.Item(8, 10).Value = 0.2637711707721372
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 17.47796633333333
$ws.Cells.Item(8, 14).Value = 52.433899
$ws.Cells.Item(8, 15).Value = 0.2471801973933361
$ws.Cells.Item(8, 16).Value = 0.2471801973933361
$ws.Cells.Item(8, 17).Value = 694.5409176071271
$ws.Cells.Item(8, 18).Value = 6250.868258464144
$ws.Cells.Item(8, 19).Value = 0.06519901005812825
$ws.Cells.Item(8, 20).Value = 0.06519901005812824

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 39.73808533333334
$ws.Cells.Item(9, 8).Value = 119.214256
$ws.Cells.Item(9, 9).Value = 0.2637711707721372
$ws.Cells.Item(9, 10).Value = 0.2637711707721372
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 12.229533
$ws.Cells.Item(9, 14).Value = 36.688599
$ws.Cells.Item(9, 15).Value = 0.17295481197965
$ws.Cells.Item(9, 16).Value = 0.17295481197965
$ws.Cells.Item(9, 17).Value = 485.9782259408161
$ws.Cells.Item(9, 18).Value = 4373.804033467344
$ws.Cells.Item(9, 19).Value = 0.04562049324654714
$ws.Cells.Item(9, 20).Value = 0.04562049324654713

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 34.33650333333333
$ws.Cells.Item(10, 8).Value = 103.00951
$ws.Cells.Item(10, 9).Value = 0.227916861330445
$ws.Cells.Item(10, 10).Value = 0.227916861330445
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 19.00842133333333
$ws.Cells.Item(10, 14).Value = 57.02526399999999
$ws.Cells.Item(10, 15).Value = 0.2688244872258518
$ws.Cells.Item(10, 16).Value = 0.2688244872258518
$ws.Cells.Item(10, 17).Value = 652.6827224734044
$ws.Cells.Item(10, 18).Value = 5874.144502260639
$ws.Cells.Item(10, 19).Value = 0.06126963337728247
$ws.Cells.Item(10, 20).Value = 0.06126963337728245

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 34.33650333333333
$ws.Cells.Item(11, 8).Value = 103.00951
$ws.Cells.Item(11, 9).Value = 0.227916861330445
$ws.Cells.Item(11, 10).Value = 0.227916861330445
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 21.993491
$ws.Cells.Item(11, 14).Value = 65.980473
$ws.Cells.Item(11, 15).Value = 0.3110405034011621
$ws.Cells.Item(11, 16).Value = 0.311040503401162
$ws.Cells.Item(11, 17).Value = 755.1795770331368
$ws.Cells.Item(11, 18).Value = 6796.61619329823
$ws.Cells.Item(11, 19).Value = 0.07089137528183447
$ws.Cells.Item(11, 20).Value = 0.07089137528183447

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 34.33650333333333
$ws.Cells.Item(12, 8).Value = 103.00951
$ws.Cells.Item(12, 9).Value = 0.227916861330445
$ws.Cells.Item(12, 10).Value = 0.227916861330445
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 17.47796633333333
$ws.Cells.Item(12, 14).Value = 52.433899
$ws.Cells.Item(12, 15).Value = 0.2471801973933361
$ws.Cells.Item(12, 16).Value = 0.2471801973933361
$ws.Cells.Item(12, 17).Value = 600.1322492643877
$ws.Cells.Item(12, 18).Value = 5401.19024337949
$ws.Cells.Item(12, 19).Value = 0.05633653477292902
$ws.Cells.Item(12, 20).Value = 0.05633653477292902

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 34.33650333333333
$ws.Cells.Item(13, 8).Value = 103.00951
$ws.Cells.Item(13, 9).Value = 0.227916861330445
$ws.Cells.Item(13, 10).Value = 0.227916861330445
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 12.229533
$ws.Cells.Item(13, 14).Value = 36.688599
$ws.Cells.Item(13, 15).Value = 0.17295481197965
$ws.Cells.Item(13, 16).Value = 0.17295481197965
$ws.Cells.Item(13, 17).Value = 419.91940061961
$ws.Cells.Item(13, 18).Value = 3779.274605576491
$ws.Cells.Item(13, 19).Value = 0.03941931789839908
$ws.Cells.Item(13, 20).Value = 0.03941931789839907

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 32.95455566666667
$ws.Cells.Item(14, 8).Value = 98.86366700000001
$ws.Cells.Item(14, 9).Value = 0.2187438488180198
$ws.Cells.Item(14, 10).Value = 0.2187438488180198
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 19.00842133333333
$ws.Cells.Item(14, 14).Value = 57.02526399999999
$ws.Cells.Item(14, 15).Value = 0.2688244872258518
$ws.Cells.Item(14, 16).Value = 0.2688244872258518
$ws.Cells.Item(14, 17).Value = 626.4140789647876
$ws.Cells.Item(14, 18).Value = 5637.726710683088
$ws.Cells.Item(14, 19).Value = 0.05880370299231343
$ws.Cells.Item(14, 20).Value = 0.0588037029923134

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 32.95455566666667
$ws.Cells.Item(15, 8).Value = 98.86366700000001
$ws.Cells.Item(15, 9).Value = 0.2187438488180198
$ws.Cells.Item(15, 10).Value = 0.2187438488180198
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 21.993491
$ws.Cells.Item(15, 14).Value = 65.980473
$ws.Cells.Item(15, 15).Value = 0.3110405034011621
$ws.Cells.Item(15, 16).Value = 0.311040503401162
$ws.Cells.Item(15, 17).Value = 724.7857234638325
$ws.Cells.Item(15, 18).Value = 6523.071511174492
$ws.Cells.Item(15, 19).Value = 0.06803819685226456
$ws.Cells.Item(15, 20).Value = 0.06803819685226455

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 32.95455566666667
$ws.Cells.Item(16, 8).Value = 98.86366700000001
$ws.Cells.Item(16, 9).Value = 0.2187438488180198
$ws.Cells.Item(16, 10).Value = 0.2187438488180198
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 17.47796633333333
$ws.Cells.Item(16, 14).Value = 52.433899
$ws.Cells.Item(16, 15).Value = 0.2471801973933361
$ws.Cells.Item(16, 16).Value = 0.2471801973933361
$ws.Cells.Item(16, 17).Value = 575.9786144719592
$ws.Cells.Item(16, 18).Value = 5183.807530247633
$ws.Cells.Item(16, 19).Value = 0.05406914772941621
$ws.Cells.Item(16, 20).Value = 0.05406914772941621

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 32.95455566666667
$ws.Cells.Item(17, 8).Value = 98.86366700000001
$ws.Cells.Item(17, 9).Value = 0.2187438488180198
$ws.Cells.Item(17, 10).Value = 0.2187438488180198
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 12.229533
$ws.Cells.Item(17, 14).Value = 36.688599
$ws.Cells.Item(17, 15).Value = 0.17295481197965
$ws.Cells.Item(17, 16).Value = 0.17295481197965
$ws.Cells.Item(17, 17).Value = 403.0188260258371
$ws.Cells.Item(17, 18).Value = 3627.169434232534
$ws.Cells.Item(17, 19).Value = 0.03783280124402559
$ws.Cells.Item(17, 20).Value = 0.03783280124402558

Write-Output "Updated rows 2-17 with recomputed NATMI values"
